$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap country labels: "Republica de Chipre" <-> "Costa de Marfil" (rows 87/88) ---
$ws.Range("A87").Value = "Costa de Marfil"
$ws.Range("A88").Value = "Republica de Chipre"

# --- Swap country labels: "Uganda" <-> "Bahamas" (rows 157/158) ---
$ws.Range("A157").Value = "Bahamas"
$ws.Range("A158").Value = "Uganda"

# --- Update "Datos actualizados" timestamp cell ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 23:22"

# --- Row 4: Estados Unidos updated stats ---
$ws.Range("B4").Value = 734846
$ws.Range("C4").Value = 25111
$ws.Range("D4").Value = 67418
$ws.Range("E4").Value = 628649
$ws.Range("F4").Value = 13536
$ws.Range("G4").Value = 1625
$ws.Range("H4").Value = 38779

# --- Row 87 (now Costa de Marfil) updated stats ---
$ws.Range("B87").Value = 801
$ws.Range("C87").Value = 113
$ws.Range("D87").Value = 239
$ws.Range("E87").Value = 554
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 8

# --- Row 88 (now Republica de Chipre) updated stats ---
$ws.Range("B88").Value = 761
$ws.Range("C88").Value = 11
$ws.Range("D88").Value = 79
$ws.Range("E88").Value = 670
$ws.Range("F88").Value = 8
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 12

# --- Row 157 (now Bahamas) updated stats ---
$ws.Range("B157").Value = 55
$ws.Range("C157").Value = 1
$ws.Range("D157").Value = 10
$ws.Range("E157").Value = 36
$ws.Range("F157").Value = 1
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 9

# --- Row 158 (now Uganda) updated stats ---
$ws.Range("B158").Value = 55
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 22
$ws.Range("E158").Value = 33
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 0
